# Edit script: add 2022-Q3 sheet and update 总计 summary sheet
$wb = $excel.ActiveWorkbook

# --- Step 1: Insert a new worksheet "2022-Q3" right before "2022-Q2" ---
$q2 = $wb.Worksheets.Item("2022-Q2")
$new = $wb.Worksheets.Add($q2)
$new.Name = "2022-Q3"

# Re-fetch source sheet reference (index shifted after Add) and copy header +
# row formatting (styles) from the "2022-Q2" sheet, which has the same shape.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Range("A1:H15").Copy($new.Range("A1:H15"))

# --- Step 2: Populate the 2022-Q3 fund-holding data ---
$new.Range("A2").Value = 0
$new.Range("B2").Value = "'006551"
$new.Range("C2").Value = "'中庚价值领航混合"
$new.Range("D2").Value = "'118.19"
$new.Range("E2").Value = "'91.86"
$new.Range("F2").Value = "'9.45"
$new.Range("G2").Value = "'11.1690"
$new.Range("H2").Value = 2
$new.Range("A3").Value = 1
$new.Range("B3").Value = "'011174"
$new.Range("C3").Value = "'中庚价值品质一年持有期混合"
$new.Range("D3").Value = "'66.33"
$new.Range("E3").Value = "'92.24"
$new.Range("F3").Value = "'9.36"
$new.Range("G3").Value = "'6.2085"
$new.Range("H3").Value = 2
$new.Range("A4").Value = 2
$new.Range("B4").Value = "'011363"
$new.Range("C4").Value = "'南方兴润价值一年持有期混合A"
$new.Range("D4").Value = "'67.01"
$new.Range("E4").Value = "'65.71"
$new.Range("F4").Value = "'2.65"
$new.Range("G4").Value = "'1.7758"
$new.Range("H4").Value = 8
$new.Range("A5").Value = 3
$new.Range("B5").Value = "'011364"
$new.Range("C5").Value = "'南方兴润价值一年持有期混合C"
$new.Range("D5").Value = "'18.36"
$new.Range("E5").Value = "'65.71"
$new.Range("F5").Value = "'2.65"
$new.Range("G5").Value = "'0.4865"
$new.Range("H5").Value = 8
$new.Range("A6").Value = 4
$new.Range("B6").Value = "'501062"
$new.Range("C6").Value = "'南方瑞合三年定期开放混合(LOF)"
$new.Range("D6").Value = "'7.16"
$new.Range("E6").Value = "'86.12"
$new.Range("F6").Value = "'4.52"
$new.Range("G6").Value = "'0.3236"
$new.Range("H6").Value = 5
$new.Range("A7").Value = 5
$new.Range("B7").Value = "'007216"
$new.Range("C7").Value = "'浙商中华预期高股息C"
$new.Range("D7").Value = "'4.40"
$new.Range("E7").Value = "'88.55"
$new.Range("F7").Value = "'6.73"
$new.Range("G7").Value = "'0.2961"
$new.Range("H7").Value = 9
$new.Range("A8").Value = 6
$new.Range("B8").Value = "'007178"
$new.Range("C8").Value = "'浙商中华预期高股息A"
$new.Range("D8").Value = "'2.59"
$new.Range("E8").Value = "'88.55"
$new.Range("F8").Value = "'6.73"
$new.Range("G8").Value = "'0.1743"
$new.Range("H8").Value = 9
$new.Range("A9").Value = 7
$new.Range("B9").Value = "'513690"
$new.Range("C9").Value = "'博时恒生港股通高股息率ETF"
$new.Range("D9").Value = "'3.05"
$new.Range("E9").Value = "'97.26"
$new.Range("F9").Value = "'3.50"
$new.Range("G9").Value = "'0.1068"
$new.Range("H9").Value = 3
$new.Range("A10").Value = 8
$new.Range("B10").Value = "'012586"
$new.Range("C10").Value = "'南方港股创新视野一年持有混合A"
$new.Range("D10").Value = "'2.16"
$new.Range("E10").Value = "'50.74"
$new.Range("F10").Value = "'4.35"
$new.Range("G10").Value = "'0.0940"
$new.Range("H10").Value = 4
$new.Range("A11").Value = 9
$new.Range("B11").Value = "'159726"
$new.Range("C11").Value = "'华夏恒生中国内地企业高股息率ETF"
$new.Range("D11").Value = "'0.84"
$new.Range("E11").Value = "'96.48"
$new.Range("F11").Value = "'3.86"
$new.Range("G11").Value = "'0.0324"
$new.Range("H11").Value = 2
$new.Range("A12").Value = 10
$new.Range("B12").Value = "'012587"
$new.Range("C12").Value = "'南方港股创新视野一年持有混合C"
$new.Range("D12").Value = "'0.19"
$new.Range("E12").Value = "'50.74"
$new.Range("F12").Value = "'4.35"
$new.Range("G12").Value = "'0.0083"
$new.Range("H12").Value = 4
$new.Range("A13").Value = 11
$new.Range("B13").Value = "'004532"
$new.Range("C13").Value = "'民生加银中证港股通高股息精选指数A"
$new.Range("D13").Value = "'0.13"
$new.Range("E13").Value = "'92.87"
$new.Range("F13").Value = "'5.21"
$new.Range("G13").Value = "'0.0068"
$new.Range("H13").Value = 2
$new.Range("A14").Value = 12
$new.Range("B14").Value = "'040021"
$new.Range("C14").Value = "'华安大中华升级股票（QDII）"
$new.Range("D14").Value = "'0.29"
$new.Range("E14").Value = "'68.57"
$new.Range("F14").Value = "'2.17"
$new.Range("G14").Value = "'0.0063"
$new.Range("H14").Value = 9
$new.Range("A15").Value = 13
$new.Range("B15").Value = "'004533"
$new.Range("C15").Value = "'民生加银中证港股通高股息精选指数C"
$new.Range("D15").Value = "'0.08"
$new.Range("E15").Value = "'92.87"
$new.Range("F15").Value = "'5.21"
$new.Range("G15").Value = "'0.0042"
$new.Range("H15").Value = 2

# --- Step 3: Update the "总计" (summary) sheet: insert the 2022-Q3 row and
# shift all subsequent rows down by one. ---
$ws1 = $wb.Worksheets.Item("总计")

# Row 9 is brand new (sheet used to only go to row 8) -- copy formatting
# (style) from row 8's A cell so A9 keeps the same bold/centered style.
$ws1.Range("A8").Copy($ws1.Range("A9"))

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 14
$ws1.Range("D2").Value = 20.69
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 15
$ws1.Range("D3").Value = 27.57
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 18
$ws1.Range("D4").Value = 16.23
$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 3.56
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 4
$ws1.Range("D6").Value = 3.54
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "2021-Q2"
$ws1.Range("C7").Value = 2
$ws1.Range("D7").Value = 0.09
$ws1.Range("A8").Value = 6
$ws1.Range("B8").Value = "2021-Q1"
$ws1.Range("C8").Value = 2
$ws1.Range("D8").Value = 4.13
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = "2020-Q4"
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 0.01
